# Update the DTT Test Hour Log worksheet with the first logged hour entry
# for the facility routes work (replacing the old "Example 1" placeholder row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Subject / Amount of hours / Date / Description
$ws.Range("A4").Value = "Start on the assignment"
$ws.Range("B4").Value = 1.1
$ws.Range("C4").Value = (Get-Date -Year 2022 -Month 6 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D4").Value = "During this hour I oriënted myself for this assignment. First I read the assignment a few times to fully understand the requirements. After I was confident enough to start I read the Readme. The steps were easy to follow and creating the database scheme took a few minutes to determine which data types to use. After the setup of the database I wanted to know how to interact with the database so I can manipulate data. This was done very easily because the start-code already has a built-in database class. As a last step I set up the routes for the facility. The /facility route now returns all facilities with their corresponding tags. You can also create update and delete the facilities but I have not yet implemented the tags part in those routes."

# The longer description now needs top-aligned, wrapped text (new cell style).
$d4 = $ws.Range("D4")
$d4.VerticalAlignment = -4160
$d4.WrapText = $true

# Move the active selection to D4 (matches the author's last cursor position).
$ws.Range("D4").Select()
